$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $value) {
    $range = $ws.Range($cellRef)
    if ($value -match '^-?\d+(\.\d+)?$') {
        # Value looks like a plain number (e.g. "2.82"); force it to stay text
        # the way the source workbook stores it, then restore the default style
        # so no stray number formatting is left behind on the cell.
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}

$updates = @(
    @("D2", "45.676.74"),
    @("E2", "  -1.80%  "),
    @("D3", "2.423.47"),
    @("E3", "  +5.71%  "),
    @("E4", "  -0.15%  "),
    @("D5", "300.76"),
    @("E5", "  -0.89%  "),
    @("D6", "97.05"),
    @("E6", "  -3.27%  "),
    @("D7", "0.567"),
    @("E7", "  +0.35%  "),
    @("E8", "  -0.06%  "),
    @("D9", "0.517"),
    @("E9", "  -0.09%  "),
    @("D10", "34.65"),
    @("E10", "  -4.81%  "),
    @("D11", "0.0795"),
    @("E11", "  +0.72%  "),
    @("D12", "7.19"),
    @("E12", "  -1.83%  "),
    @("E13", "  +0.97%  "),
    @("D14", "2.787.72"),
    @("E14", "  +5.52%  "),
    @("D15", "2.432.46"),
    @("E15", "  +6.19%  "),
    @("D16", "14.38"),
    @("E16", "  +4.51%  "),
    @("D17", "0.843"),
    @("E17", "  +4.50%  "),
    @("D18", "45.670.93"),
    @("E18", "  -1.82%  "),
    @("D19", "13.29"),
    @("E19", "  +2.60%  "),
    @("D20", "0.0₃0954"),
    @("E20", "  +2.47%  "),
    @("D21", "6.22"),
    @("E21", "  +4.21%  "),
    @("D22", "67.56"),
    @("E22", "  +2.07%  "),
    @("D23", "244.32"),
    @("E23", "  -1.36%  "),
    @("D24", "2.82"),
    @("E24", "  -2.27%  "),
    @("D25", "1.96"),
    @("E25", "  +1.96%  "),
    @("E26", "  -0.07%  "),
    @("D27", "38.88"),
    @("E27", "  -8.44%  "),
    @("E28", "  -1.53%  "),
    @("D29", "9.83"),
    @("E29", "  +1.01%  "),
    @("D30", "3.86"),
    @("E30", "  +18.11%  "),
    @("D31", "21.47"),
    @("E31", "  +7.74%  "),
    @("D32", "5.60"),
    @("E32", "  -0.45%  "),
    @("E33", "  -1.60%  "),
    @("D34", "148.25"),
    @("E34", "  +1.10%  "),
    @("D35", "0.0780"),
    @("E35", "  -1.46%  "),
    @("D36", "2.01"),
    @("E36", "  +13.74%  "),
    @("E37", "  -0.03%  "),
    @("E38", "  -0.89%  "),
    @("D39", "15.25"),
    @("E39", "  -3.62%  "),
    @("D40", "3.90"),
    @("E40", "  -2.55%  "),
    @("D41", "0.0302"),
    @("E41", "  +0.34%  "),
    @("D42", "3.29"),
    @("E42", "  -1.34%  "),
    @("D43", "1.954.35"),
    @("E43", "  +7.58%  "),
    @("E44", "  -0.07%  "),
    @("D45", "91.28"),
    @("E45", "  +4.83%  "),
    @("E46", "  -9.26%  "),
    @("D47", "8.69"),
    @("E47", "  +10.80%  "),
    @("B48", "EnergySwap"),
    @("C48", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"),
    @("D48", "15.39"),
    @("E48", "  +15.25%  "),
    @("B49", "Aave"),
    @("C49", "https://coinranking.com/coin/ixgUfzmLR+aave-aave"),
    @("D49", "101.57"),
    @("E49", "  +6.52%  "),
    @("D50", "0.187"),
    @("E50", "  -3.23%  "),
    @("B51", "RocketPoolETH"),
    @("C51", "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"),
    @("D51", "2.659.61"),
    @("E51", "  +5.61%  ")
)

foreach ($u in $updates) {
    Set-TextValue $ws $u[0] $u[1]
}
